# Add two new roster rows (Thomas / Pilon, both "FF") to the draft board,
# matching the layout/formatting already used by the existing rows above
# them (row 82, "Frei").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last populated row (82) onto the two blank
# rows (83, 84) right below it, so the new cells pick up the same styles
# (name column border/font, rank column font, boolean-column font) as the
# rest of the table.
$ws.Range("A82:K82").Copy()
$ws.Range("A83:K83").PasteSpecial(-4122)
$ws.Range("A84:K84").PasteSpecial(-4122)

# Row 83: Thomas, rank FF, no qualifications checked.
$ws.Range("A83").Value = "Thomas"
$ws.Range("B83").Value = "FF"
$ws.Range("C83:K83").Value = $false

# Row 84: Pilon, rank FF, no qualifications checked.
$ws.Range("A84").Value = "Pilon"
$ws.Range("B84").Value = "FF"
$ws.Range("C84:K84").Value = $false
